$d = $word.ActiveDocument
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Locate the "Use case: ... medical doctors and researchers)" bullet paragraph
$found = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Use case:*researchers)*") {
        $anchor = $p
        $found = $true
    }
}

if (-not $found) {
    throw "Could not locate anchor paragraph"
}

# --- Insert first new bullet paragraph right after the anchor paragraph ---
$rng = $anchor.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$anchorIndex = $anchor.Index
$p1 = $d.Paragraphs.Item($anchorIndex + 1)
$p1.Range.InsertBefore("placeholder")
$r1 = $d.Range($p1.Range.Start, $p1.Range.End)
$xml1 = "<w:p $w><w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""1""/></w:numPr></w:pPr><w:r><w:t xml:space=""preserve"">How </w:t></w:r><w:proofErr w:type=""spellStart""/><w:r><w:t>iMedBot</w:t></w:r><w:proofErr w:type=""spellEnd""/><w:r><w:t xml:space=""preserve""> was deployed from </w:t></w:r><w:proofErr w:type=""spellStart""/><w:r><w:t>github</w:t></w:r><w:proofErr w:type=""spellEnd""/><w:r><w:t xml:space=""preserve""> to AWA.</w:t></w:r></w:p>"
$r1.InsertXML($xml1)

# --- Insert second new bullet paragraph right after the first new paragraph ---
$p1 = $d.Paragraphs.Item($anchorIndex + 1)
$rng2 = $p1.Range
$rng2.Collapse(0)
$rng2.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($anchorIndex + 2)
$p2.Range.InsertBefore("placeholder")
$r2 = $d.Range($p2.Range.Start, $p2.Range.End)
$xml2 = "<w:p $w><w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""1""/></w:numPr></w:pPr><w:r><w:t>How to access iMedBot.</w:t></w:r></w:p>"
$r2.InsertXML($xml2)
